# "Vision de la aplicacion.docx" - "Added MongoDB Connection and Videos WIP"
#
# 1) Fix the "Contrasena" -> "Contrase(n~)a" typo (the auto spell-check
#    proofErr markers that wrapped the misspelled word are cleared once
#    the word is corrected).
# 2) The hidden "_GoBack" last-edit bookmark moves from the
#    "Suscripciones: Nombre" bullet to the "Like" bullet, since that is
#    where the final edit in this session happens (the WIP video fields).
# 3) The placeholder "Ruta video" / "Ruta Miniatura" bullets under
#    "Videos:" are removed (still a work in progress).

$d = $word.ActiveDocument
$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

function Get-ParaIndexByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $t = $doc.Paragraphs.Item($i).Range.Text
        $t = $t.Replace([char]13, "")
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

function Set-ParaXml($doc, $index, $innerXml) {
    $p = $doc.Paragraphs.Item($index)
    $full = $p.Range
    $xml = '<w:p xmlns:w="' + $wNs + '">' + $innerXml + '</w:p>'
    $full.InsertXML($xml)
}

$pPrLvl1 = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr>'
$pPrLvl2 = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr>'

# 1. "Contrasena" -> "Contrasenia" (n with tilde), drop the proofErr wrap.
$idx = Get-ParaIndexByText $d "Contrasena"
$inner = $pPrLvl1 + '<w:r><w:t>Contrase&#241;a</w:t></w:r>'
Set-ParaXml $d $idx $inner

# 2. Drop the "_GoBack" bookmark from "Suscripciones: Nombre".
$idx = Get-ParaIndexByText $d "Suscripciones: Nombre"
$inner = $pPrLvl1 + '<w:r><w:t>Suscripciones: Nombre</w:t></w:r>'
Set-ParaXml $d $idx $inner

# 3. Move the "_GoBack" bookmark onto the "Like" bullet (keep its
#    existing proofErr spell-check wrap around "Like").
$idx = Get-ParaIndexByText $d "Like"
$inner = $pPrLvl2 + '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:proofErr w:type="spellStart"/><w:r><w:t>Like</w:t></w:r><w:proofErr w:type="spellEnd"/>'
Set-ParaXml $d $idx $inner

# 4. Remove the still-unfinished "Ruta video" / "Ruta Miniatura" bullets.
$idx = Get-ParaIndexByText $d "Ruta video"
if ($idx -gt 0) { $d.Paragraphs.Item($idx).Range.Delete() }

$idx = Get-ParaIndexByText $d "Ruta Miniatura"
if ($idx -gt 0) { $d.Paragraphs.Item($idx).Range.Delete() }
